$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.004.86"
$ws.Range("E2").Value = "  -4.34%  "

$ws.Range("D3").Value = "2.939.19"
$ws.Range("E3").Value = "  -7.13%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "568.07"
$ws.Range("E5").Value = "  -3.73%  "

$ws.Range("D6").Value = "121.76"
$ws.Range("E6").Value = "  -9.51%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "2.935.45"
$ws.Range("E8").Value = "  -7.18%  "

$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -4.26%  "

$ws.Range("D10").Value = "0.130"
$ws.Range("E10").Value = "  -6.97%  "

$ws.Range("D11").Value = "5.03"
$ws.Range("E11").Value = "  -3.68%  "

$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  -4.46%  "

$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -6.97%  "

$ws.Range("D14").Value = "32.07"
$ws.Range("E14").Value = "  -7.57%  "

$ws.Range("E15").Value = "  -0.93%  "

$ws.Range("D16").Value = "3.429.59"
$ws.Range("E16").Value = "  -7.03%  "

$ws.Range("D17").Value = "60.007.02"
$ws.Range("E17").Value = "  -4.34%  "

$ws.Range("D18").Value = "2.944.88"
$ws.Range("E18").Value = "  -6.97%  "

$ws.Range("D19").Value = "6.09"
$ws.Range("E19").Value = "  -6.83%  "

$ws.Range("D20").Value = "424.30"
$ws.Range("E20").Value = "  -7.50%  "

$ws.Range("D21").Value = "12.88"
$ws.Range("E21").Value = "  -7.91%  "

$ws.Range("D22").Value = "0.652"
$ws.Range("E22").Value = "  -5.70%  "

$ws.Range("D23").Value = "6.91"
$ws.Range("E23").Value = "  -9.15%  "

$ws.Range("D24").Value = "12.69"
$ws.Range("E24").Value = "  -4.97%  "

$ws.Range("D25").Value = "77.83"
$ws.Range("E25").Value = "  -5.95%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "2.49"
$ws.Range("E28").Value = "  -6.33%  "

$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  -8.30%  "

$ws.Range("D30").Value = "1.85"
$ws.Range("E30").Value = "  -8.25%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "24.92"
$ws.Range("E31").Value = "  -7.88%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "6.01"
$ws.Range("E32").Value = "  -10.96%  "

$ws.Range("D33").Value = "0.0912"
$ws.Range("E33").Value = "  -10.63%  "

$ws.Range("D34").Value = "2.15"
$ws.Range("E34").Value = "  -8.67%  "

$ws.Range("E35").Value = "  -9.13%  "

$ws.Range("D36").Value = "5.48"
$ws.Range("E36").Value = "  -5.23%  "

$ws.Range("D37").Value = "49.19"
$ws.Range("E37").Value = "  -3.97%  "

$ws.Range("D38").Value = "0.0₃0638"
$ws.Range("E38").Value = "  -8.90%  "

$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").Value = "7.77"
$ws.Range("E39").Value = "  -3.30%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0352"
$ws.Range("E40").Value = "  -8.84%  "

$ws.Range("D41").Value = "0.107"
$ws.Range("E41").Value = "  -3.67%  "

$ws.Range("D42").Value = "372.36"
$ws.Range("E42").Value = "  -7.57%  "

$ws.Range("D43").Value = "2.602.38"
$ws.Range("E43").Value = "  -6.60%  "

$ws.Range("D44").Value = "2.36"
$ws.Range("E44").Value = "  -9.59%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "0.232"
$ws.Range("E46").Value = "  -7.20%  "

$ws.Range("D47").Value = "118.61"
$ws.Range("E47").Value = "  -4.73%  "

$ws.Range("E48").Value = "  -8.05%  "

$ws.Range("E49").Value = "  -5.39%  "

$ws.Range("D50").Value = "22.94"
$ws.Range("E50").Value = "  -8.73%  "

$ws.Range("D51").Value = "30.56"
$ws.Range("E51").Value = "  -10.74%  "
